# Apply the "technology types removed and duration revised" edit:
# revise the duration values (column C) for several technology rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C15").Value = 0.3
$ws.Range("C18").Value = 0.3
$ws.Range("C19").Value = 0.2
$ws.Range("C20").Value = 0.2
$ws.Range("C21").Value = 0.3
$ws.Range("C22").Value = 0.2
$ws.Range("C27").Value = 0.3
$ws.Range("C28").Value = 0.3
$ws.Range("C29").Value = 0.3

# Update the view state (last-saved scroll position / zoom / selection)
# to match the author's cursor position when they saved the revision.
[void]$ws.Range("C39").Select()
$excel.ActiveWindow.Zoom = 100
